$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the existing hyperlink on A3 (testnum123@yopmail.com) before re-writing it.
$ws.Range("A3").Hyperlinks.Delete()

# Write the new cell values (row 4 is a new data row copied/adapted from row 3,
# row 3's email is renumbered from testnum123 to testnum12).
$ws.Range("A4").Value = "testnum456@yopmail.com"
$ws.Range("B4").Value = "Jis"
$ws.Range("C4").Value = "Jay"
$ws.Range("A3").Value = "testnum12@yopmail.com"
$ws.Range("D4").Value = "BT@2021"

# Re-create the hyperlinks (mailto: links) and restore the "Hyperlink" cell style.
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:testnum12@yopmail.com")
$ws.Range("A3").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:testnum456@yopmail.com")
$ws.Range("A4").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:BT@2021")
$ws.Range("D4").Style = "Hyperlink"

# Move the active selection to B13, matching the saved view state.
$ws.Range("B13").Select()
